$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2089356110381078
$ws.Range("C2").Value = 0.5308804204993429
$ws.Range("J2").Value = 0.007884362680683311
$ws.Range("P2").Value = 0.1406044678055191
$ws.Range("S2").Value = 0.1116951379763469
$ws.Range("B3").Value = 0.002375296912114014
$ws.Range("C3").Value = 0.02850356294536817
$ws.Range("J3").Value = 0.05225653206650831
$ws.Range("P3").Value = 0.7125890736342043
$ws.Range("S3").Value = 0.2042755344418052
$ws.Range("J4").Value = 0.05737704918032787
$ws.Range("P4").Value = 0.680327868852459
$ws.Range("S4").Value = 0.2622950819672131
$ws.Range("B6").Value = 0.06396255850234009
$ws.Range("D6").Value = 0.0109204368174727
$ws.Range("F6").Value = 0.0655226209048362
$ws.Range("J6").Value = 0.2636505460218408
$ws.Range("O6").Value = 0.031201248049922
$ws.Range("Q6").Value = 0.1669266770670827
$ws.Range("R6").Value = 0.0686427457098284
$ws.Range("S6").Value = 0.3291731669266771
$ws.Range("B7").Value = 0.09677419354838709
$ws.Range("D7").Value = 0.02656546489563567
$ws.Range("F7").Value = 0.06072106261859583
$ws.Range("J7").Value = 0.1518026565464896
$ws.Range("O7").Value = 0.01138519924098672
$ws.Range("Q7").Value = 0.1480075901328273
$ws.Range("R7").Value = 0.08159392789373814
$ws.Range("S7").Value = 0.4231499051233397
$ws.Range("B8").Value = 0.09334485738980121
$ws.Range("D8").Value = 0.0233362143474503
$ws.Range("E8").Value = 0.001728608470181504
$ws.Range("F8").Value = 0.05963699222126188
$ws.Range("J8").Value = 0.118409680207433
$ws.Range("O8").Value = 0.01901469317199654
$ws.Range("Q8").Value = 0.1624891961970614
$ws.Range("R8").Value = 0.09075194468452895
$ws.Range("S8").Value = 0.4312878133102852
$ws.Range("B9").Value = 0.09073724007561437
$ws.Range("D9").Value = 0.01701323251417769
$ws.Range("F9").Value = 0.06238185255198488
$ws.Range("J9").Value = 0.1153119092627599
$ws.Range("O9").Value = 0.02268431001890359
$ws.Range("Q9").Value = 0.1965973534971645
$ws.Range("R9").Value = 0.06994328922495274
$ws.Range("S9").Value = 0.4253308128544424
$ws.Range("B10").Value = 0.1089542892924233
$ws.Range("D10").Value = 0.02128991859737007
$ws.Range("E10").Value = 0.000939261114589856
$ws.Range("F10").Value = 0.08015028177833437
$ws.Range("J10").Value = 0.1102066374452098
$ws.Range("O10").Value = 0.01408891671884784
$ws.Range("Q10").Value = 0.2038196618659988
$ws.Range("R10").Value = 0.08547276142767689
$ws.Range("S10").Value = 0.3750782717595492
$ws.Range("G11").Value = 0.1124031007751938
$ws.Range("J11").Value = 0.103359173126615
$ws.Range("K11").Value = 0.1640826873385013
$ws.Range("L11").Value = 0.6072351421188631
$ws.Range("S11").Value = 0.01291989664082687
$ws.Range("G12").Value = 0.7525562372188139
$ws.Range("J12").Value = 0.1697341513292434
$ws.Range("K12").Value = 0.01431492842535787
$ws.Range("L12").Value = 0.03067484662576687
$ws.Range("S12").Value = 0.032719836400818
$ws.Range("G13").Value = 0.6615384615384615
$ws.Range("J13").Value = 0.2538461538461538
$ws.Range("S13").Value = 0.08461538461538462
$ws.Range("F15").Value = 0.01211072664359862
$ws.Range("H15").Value = 0.1608996539792387
$ws.Range("I15").Value = 0.06920415224913495
$ws.Range("J15").Value = 0.3892733564013841
$ws.Range("K15").Value = 0.05709342560553633
$ws.Range("M15").Value = 0.008650519031141869
$ws.Range("N15").Value = 0.001730103806228374
$ws.Range("O15").Value = 0.06747404844290658
$ws.Range("S15").Value = 0.2335640138408304
$ws.Range("F16").Value = 0.03131524008350731
$ws.Range("H16").Value = 0.1920668058455115
$ws.Range("I16").Value = 0.07515657620041753
$ws.Range("J16").Value = 0.3716075156576201
$ws.Range("K16").Value = 0.1294363256784969
$ws.Range("M16").Value = 0.02505219206680585
$ws.Range("N16").Value = 0.00208768267223382
$ws.Range("O16").Value = 0.05010438413361169
$ws.Range("S16").Value = 0.1231732776617954
$ws.Range("F17").Value = 0.02044444444444445
$ws.Range("H17").Value = 0.1857777777777778
$ws.Range("I17").Value = 0.1102222222222222
$ws.Range("J17").Value = 0.3768888888888889
$ws.Range("K17").Value = 0.104
$ws.Range("M17").Value = 0.02044444444444445
$ws.Range("O17").Value = 0.07555555555555556
$ws.Range("S17").Value = 0.1066666666666667
$ws.Range("F18").Value = 0.02589641434262948
$ws.Range("H18").Value = 0.1673306772908366
$ws.Range("I18").Value = 0.1035856573705179
$ws.Range("J18").Value = 0.4183266932270917
$ws.Range("K18").Value = 0.1055776892430279
$ws.Range("M18").Value = 0.02589641434262948
$ws.Range("N18").Value = 0.00199203187250996
$ws.Range("O18").Value = 0.06772908366533864
$ws.Range("S18").Value = 0.08366533864541832
$ws.Range("F19").Value = 0.02301319423135931
$ws.Range("H19").Value = 0.2101871739797484
$ws.Range("I19").Value = 0.08560908254065665
$ws.Range("J19").Value = 0.360846885547714
$ws.Range("K19").Value = 0.1107701749002762
$ws.Range("M19").Value = 0.02516109235961951
$ws.Range("N19").Value = 0.001841055538508745
$ws.Range("O19").Value = 0.07241485118134397
$ws.Range("S19").Value = 0.1101564897207732
